# KIBON-2068: FR Uebersetzungen + neue Titlen fuer die die waren noch hardcodiert
#
# The "Data" sheet had three hardcoded German header labels (for the
# monthly-share columns) that should instead be placeholder tokens like
# all the other report headers, so the reporting engine can localize them
# (e.g. into French).
#
# O8: "Tage Monat"        -> {tageMonatTitle}
# P8: "Tage Intervall"    -> {tageIntervallTitle}
# Q8: "Anteil des Monats" -> {anteilMonatKantonTitle}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("O8").Value = "{tageMonatTitle}"
$ws.Range("P8").Value = "{tageIntervallTitle}"
$ws.Range("Q8").Value = "{anteilMonatKantonTitle}"

# Restore the selection to the merged Q8:Q9 cell, matching the state the
# workbook was left in when it was last saved.
$ws.Range("Q8:Q9").Select()
